$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New daily rows for 02-10-2021 .. 06-10-2021, continuing the existing series
# with the same values as the rest of the October block.
$xlPasteValues = -4163

$newRows = @(
    @{ Date = "02-10-2021"; B = 1694; C = 2114; D = 12227; E = 2134; F = 3696; G = 7100 },
    @{ Date = "03-10-2021"; B = 1694; C = 2114; D = 12227; E = 2134; F = 3696; G = 7100 },
    @{ Date = "04-10-2021"; B = 1694; C = 2114; D = 12227; E = 2134; F = 3696; G = 7100 },
    @{ Date = "05-10-2021"; B = 1694; C = 2114; D = 12227; E = 2134; F = 3696; G = 7100 },
    @{ Date = "06-10-2021"; B = 1694; C = 2114; D = 12227; E = 2134; F = 3696; G = 7100 }
)

$startRow = 276
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Column A holds a dd-mm-yyyy text label. Writing the literal string via
    # .Value/.Formula would get auto-recognized as a date (since both the
    # day and month are <= 12) and reformatted as a number. Instead, enter
    # it as a text formula and convert it to a plain value via
    # Copy/PasteSpecial(xlPasteValues) so it lands as plain text, matching
    # how the rest of the column is stored.
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.Formula = '="' + $row.Date + '"'
    $dateCell.Copy() | Out-Null
    $dateCell.PasteSpecial($xlPasteValues) | Out-Null

    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
}

$excel.CutCopyMode = $false
